# "Replaced broken assets from scene 1"
#
# The feedback note for row 9 ("8. All models successfully drawn at proper
# location, orientation & scale") changes from a "partially complete /
# some submeshes aren't drawing" remark to "Complete" now that the broken
# scene-1 assets have been replaced.
#
# A couple of other feedback notes on the rubric sheet are also touched as
# part of the same pass: the note for row 14 moves from column E to column D,
# a new "d" note is added for row 30, and "Complete" notes are added for
# rows 46 and 47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 ("8. All models successfully drawn..."): feedback updated to Complete.
$ws.Range("D9").Value = "Complete"

# Row 14 ("Support for Hot-Swapping multiple levels..."): feedback note moved
# from column E to column D (still "Complete").
$ws.Range("E14").ClearContents()
$ws.Range("D14").Value = "Complete"

# Row 30 ("Effectively use a new type of Shader."): new short note "d" added.
$ws.Range("E30").Value = "d"

# Row 46 ("Play some Level Music...") and Row 47 ("Manage your project in a
# private GIT repo..."): new "Complete" feedback notes added.
$ws.Range("D46").Value = "Complete"
$ws.Range("D47").Value = "Complete"

# Selection moved to A10 (previously H11, with the view scrolled so A10 was
# pinned at the top-left).
$ws.Range("A10").Select()
